# Apply commit "Update latest output (run 59)" to optimisation_result.xlsx
# Sheet1 = "Schedule" (pump run schedule), Sheet2 = "Detailed" (price / pump-status series).
$wb = $excel.ActiveWorkbook

# ============================================================
# Schedule sheet: row 2 values replaced, rows 3-5 newly added
# ============================================================
$schedule = $wb.Worksheets.Item("Schedule")

$schedule.Cells.Item(2, 1).Value = 46039.02083333334
$schedule.Cells.Item(2, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$schedule.Cells.Item(2, 2).Value = 46039.1875
$schedule.Cells.Item(2, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$schedule.Cells.Item(2, 3).Value = 4
$schedule.Cells.Item(2, 4).Value = 15.12
$schedule.Cells.Item(2, 5).Value = 347.06108775
$schedule.Cells.Item(2, 6).Value = 22.95377564484127

$schedule.Cells.Item(3, 1).Value = 46039.29166666666
$schedule.Cells.Item(3, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$schedule.Cells.Item(3, 2).Value = 46039.5
$schedule.Cells.Item(3, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$schedule.Cells.Item(3, 3).Value = 5
$schedule.Cells.Item(3, 4).Value = 18.9
$schedule.Cells.Item(3, 5).Value = 313.99407
$schedule.Cells.Item(3, 6).Value = 16.61344285714286

$schedule.Cells.Item(4, 1).Value = 46039.52083333334
$schedule.Cells.Item(4, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$schedule.Cells.Item(4, 2).Value = 46039.8125
$schedule.Cells.Item(4, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$schedule.Cells.Item(4, 3).Value = 7
$schedule.Cells.Item(4, 4).Value = 26.46
$schedule.Cells.Item(4, 5).Value = 9.565959000000007
$schedule.Cells.Item(4, 6).Value = 0.3615252834467123

$schedule.Cells.Item(5, 1).Value = 46040.29166666666
$schedule.Cells.Item(5, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$schedule.Cells.Item(5, 2).Value = 46040.79166666666
$schedule.Cells.Item(5, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$schedule.Cells.Item(5, 3).Value = 12
$schedule.Cells.Item(5, 4).Value = 45.36
$schedule.Cells.Item(5, 5).Value = -95.68104975
$schedule.Cells.Item(5, 6).Value = -2.10937058531746

# ============================================================
# Detailed sheet: targeted cell corrections for existing rows 2-49
# ============================================================
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Cells.Item(3, 5).Value = "ON"
$detailed.Cells.Item(4, 5).Value = "ON"
$detailed.Cells.Item(5, 5).Value = "ON"
$detailed.Cells.Item(6, 5).Value = "ON"
$detailed.Cells.Item(7, 5).Value = "ON"
$detailed.Cells.Item(8, 5).Value = "ON"
$detailed.Cells.Item(9, 5).Value = "ON"
$detailed.Cells.Item(10, 5).Value = "ON"
$detailed.Cells.Item(15, 5).Value = "OFF"
$detailed.Cells.Item(25, 2).Value = 36.06046
$detailed.Cells.Item(26, 2).Value = 36.06092
$detailed.Cells.Item(26, 5).Value = "OFF"
$detailed.Cells.Item(27, 2).Value = 36.06045
$detailed.Cells.Item(27, 3).Value = "historical"
$detailed.Cells.Item(28, 2).Value = -9.99
$detailed.Cells.Item(28, 3).Value = "historical"
$detailed.Cells.Item(29, 2).Value = -6.8
$detailed.Cells.Item(29, 3).Value = "historical"
$detailed.Cells.Item(30, 2).Value = -9.99
$detailed.Cells.Item(30, 3).Value = "historical"
$detailed.Cells.Item(31, 2).Value = -7
$detailed.Cells.Item(32, 2).Value = -5.51011
$detailed.Cells.Item(33, 2).Value = -5.01
$detailed.Cells.Item(34, 2).Value = 0.7
$detailed.Cells.Item(35, 2).Value = 3.23527
$detailed.Cells.Item(36, 2).Value = 0
$detailed.Cells.Item(37, 2).Value = -3.09373
$detailed.Cells.Item(38, 2).Value = 0.00013
$detailed.Cells.Item(39, 2).Value = 0
$detailed.Cells.Item(40, 2).Value = 17.20923
$detailed.Cells.Item(41, 2).Value = 57.16514
$detailed.Cells.Item(41, 5).Value = "OFF"
$detailed.Cells.Item(42, 2).Value = 57.18142
$detailed.Cells.Item(42, 5).Value = "OFF"
$detailed.Cells.Item(43, 2).Value = 57.16514
$detailed.Cells.Item(43, 5).Value = "OFF"
$detailed.Cells.Item(44, 2).Value = 46.73943
$detailed.Cells.Item(44, 5).Value = "OFF"
$detailed.Cells.Item(45, 2).Value = 46.80565
$detailed.Cells.Item(45, 5).Value = "OFF"
$detailed.Cells.Item(46, 2).Value = 36.05843
$detailed.Cells.Item(46, 5).Value = "OFF"
$detailed.Cells.Item(47, 2).Value = 47.35317

# ============================================================
# Detailed sheet: brand new rows 50-97 (forecast data for 46040 = next day)
# ============================================================
$detailed.Cells.Item(50, 1).Value = 46040
$detailed.Cells.Item(50, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(50, 2).Value = 47.81285
$detailed.Cells.Item(50, 3).Value = "forecast"
$detailed.Cells.Item(50, 4).Value = 46040
$detailed.Cells.Item(50, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(50, 5).Value = "OFF"

$detailed.Cells.Item(51, 1).Value = 46040.02083333334
$detailed.Cells.Item(51, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(51, 2).Value = 36.06032
$detailed.Cells.Item(51, 3).Value = "forecast"
$detailed.Cells.Item(51, 4).Value = 46040
$detailed.Cells.Item(51, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(51, 5).Value = "OFF"

$detailed.Cells.Item(52, 1).Value = 46040.04166666666
$detailed.Cells.Item(52, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(52, 2).Value = 56.98
$detailed.Cells.Item(52, 3).Value = "forecast"
$detailed.Cells.Item(52, 4).Value = 46040
$detailed.Cells.Item(52, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(52, 5).Value = "OFF"

$detailed.Cells.Item(53, 1).Value = 46040.0625
$detailed.Cells.Item(53, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(53, 2).Value = 56.98
$detailed.Cells.Item(53, 3).Value = "forecast"
$detailed.Cells.Item(53, 4).Value = 46040
$detailed.Cells.Item(53, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(53, 5).Value = "OFF"

$detailed.Cells.Item(54, 1).Value = 46040.08333333334
$detailed.Cells.Item(54, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(54, 2).Value = 56.97996
$detailed.Cells.Item(54, 3).Value = "forecast"
$detailed.Cells.Item(54, 4).Value = 46040
$detailed.Cells.Item(54, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(54, 5).Value = "OFF"

$detailed.Cells.Item(55, 1).Value = 46040.10416666666
$detailed.Cells.Item(55, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(55, 2).Value = 57.0602
$detailed.Cells.Item(55, 3).Value = "forecast"
$detailed.Cells.Item(55, 4).Value = 46040
$detailed.Cells.Item(55, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(55, 5).Value = "OFF"

$detailed.Cells.Item(56, 1).Value = 46040.125
$detailed.Cells.Item(56, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(56, 2).Value = 57.06021
$detailed.Cells.Item(56, 3).Value = "forecast"
$detailed.Cells.Item(56, 4).Value = 46040
$detailed.Cells.Item(56, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(56, 5).Value = "OFF"

$detailed.Cells.Item(57, 1).Value = 46040.14583333334
$detailed.Cells.Item(57, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(57, 2).Value = 57.06022
$detailed.Cells.Item(57, 3).Value = "forecast"
$detailed.Cells.Item(57, 4).Value = 46040
$detailed.Cells.Item(57, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(57, 5).Value = "OFF"

$detailed.Cells.Item(58, 1).Value = 46040.16666666666
$detailed.Cells.Item(58, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(58, 2).Value = 57.06022
$detailed.Cells.Item(58, 3).Value = "forecast"
$detailed.Cells.Item(58, 4).Value = 46040
$detailed.Cells.Item(58, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(58, 5).Value = "OFF"

$detailed.Cells.Item(59, 1).Value = 46040.1875
$detailed.Cells.Item(59, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(59, 2).Value = 57.06003
$detailed.Cells.Item(59, 3).Value = "forecast"
$detailed.Cells.Item(59, 4).Value = 46040
$detailed.Cells.Item(59, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(59, 5).Value = "OFF"

$detailed.Cells.Item(60, 1).Value = 46040.20833333334
$detailed.Cells.Item(60, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(60, 2).Value = 57.06003
$detailed.Cells.Item(60, 3).Value = "forecast"
$detailed.Cells.Item(60, 4).Value = 46040
$detailed.Cells.Item(60, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(60, 5).Value = "OFF"

$detailed.Cells.Item(61, 1).Value = 46040.22916666666
$detailed.Cells.Item(61, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(61, 2).Value = 57.1
$detailed.Cells.Item(61, 3).Value = "forecast"
$detailed.Cells.Item(61, 4).Value = 46040
$detailed.Cells.Item(61, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(61, 5).Value = "OFF"

$detailed.Cells.Item(62, 1).Value = 46040.25
$detailed.Cells.Item(62, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(62, 2).Value = 57.06003
$detailed.Cells.Item(62, 3).Value = "forecast"
$detailed.Cells.Item(62, 4).Value = 46040
$detailed.Cells.Item(62, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(62, 5).Value = "OFF"

$detailed.Cells.Item(63, 1).Value = 46040.27083333334
$detailed.Cells.Item(63, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(63, 2).Value = 36.06
$detailed.Cells.Item(63, 3).Value = "forecast"
$detailed.Cells.Item(63, 4).Value = 46040
$detailed.Cells.Item(63, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(63, 5).Value = "OFF"

$detailed.Cells.Item(64, 1).Value = 46040.29166666666
$detailed.Cells.Item(64, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(64, 2).Value = 24.71692
$detailed.Cells.Item(64, 3).Value = "forecast"
$detailed.Cells.Item(64, 4).Value = 46040
$detailed.Cells.Item(64, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(64, 5).Value = "ON"

$detailed.Cells.Item(65, 1).Value = 46040.3125
$detailed.Cells.Item(65, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(65, 2).Value = 0.51
$detailed.Cells.Item(65, 3).Value = "forecast"
$detailed.Cells.Item(65, 4).Value = 46040
$detailed.Cells.Item(65, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(65, 5).Value = "ON"

$detailed.Cells.Item(66, 1).Value = 46040.33333333334
$detailed.Cells.Item(66, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(66, 2).Value = -4.58324
$detailed.Cells.Item(66, 3).Value = "forecast"
$detailed.Cells.Item(66, 4).Value = 46040
$detailed.Cells.Item(66, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(66, 5).Value = "ON"

$detailed.Cells.Item(67, 1).Value = 46040.35416666666
$detailed.Cells.Item(67, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(67, 2).Value = 0.62395
$detailed.Cells.Item(67, 3).Value = "forecast"
$detailed.Cells.Item(67, 4).Value = 46040
$detailed.Cells.Item(67, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(67, 5).Value = "ON"

$detailed.Cells.Item(68, 1).Value = 46040.375
$detailed.Cells.Item(68, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(68, 2).Value = 2.19121
$detailed.Cells.Item(68, 3).Value = "forecast"
$detailed.Cells.Item(68, 4).Value = 46040
$detailed.Cells.Item(68, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(68, 5).Value = "ON"

$detailed.Cells.Item(69, 1).Value = 46040.39583333334
$detailed.Cells.Item(69, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(69, 2).Value = -0.9231200000000001
$detailed.Cells.Item(69, 3).Value = "forecast"
$detailed.Cells.Item(69, 4).Value = 46040
$detailed.Cells.Item(69, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(69, 5).Value = "ON"

$detailed.Cells.Item(70, 1).Value = 46040.41666666666
$detailed.Cells.Item(70, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(70, 2).Value = -0.9403
$detailed.Cells.Item(70, 3).Value = "forecast"
$detailed.Cells.Item(70, 4).Value = 46040
$detailed.Cells.Item(70, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(70, 5).Value = "ON"

$detailed.Cells.Item(71, 1).Value = 46040.4375
$detailed.Cells.Item(71, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(71, 2).Value = -5.4861
$detailed.Cells.Item(71, 3).Value = "forecast"
$detailed.Cells.Item(71, 4).Value = 46040
$detailed.Cells.Item(71, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(71, 5).Value = "ON"

$detailed.Cells.Item(72, 1).Value = 46040.45833333334
$detailed.Cells.Item(72, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(72, 2).Value = -5.51
$detailed.Cells.Item(72, 3).Value = "forecast"
$detailed.Cells.Item(72, 4).Value = 46040
$detailed.Cells.Item(72, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(72, 5).Value = "ON"

$detailed.Cells.Item(73, 1).Value = 46040.47916666666
$detailed.Cells.Item(73, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(73, 2).Value = -5.77472
$detailed.Cells.Item(73, 3).Value = "forecast"
$detailed.Cells.Item(73, 4).Value = 46040
$detailed.Cells.Item(73, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(73, 5).Value = "ON"

$detailed.Cells.Item(74, 1).Value = 46040.5
$detailed.Cells.Item(74, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(74, 2).Value = -5.63201
$detailed.Cells.Item(74, 3).Value = "forecast"
$detailed.Cells.Item(74, 4).Value = 46040
$detailed.Cells.Item(74, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(74, 5).Value = "ON"

$detailed.Cells.Item(75, 1).Value = 46040.52083333334
$detailed.Cells.Item(75, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(75, 2).Value = -6.30444
$detailed.Cells.Item(75, 3).Value = "forecast"
$detailed.Cells.Item(75, 4).Value = 46040
$detailed.Cells.Item(75, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(75, 5).Value = "ON"

$detailed.Cells.Item(76, 1).Value = 46040.54166666666
$detailed.Cells.Item(76, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(76, 2).Value = -7
$detailed.Cells.Item(76, 3).Value = "forecast"
$detailed.Cells.Item(76, 4).Value = 46040
$detailed.Cells.Item(76, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(76, 5).Value = "ON"

$detailed.Cells.Item(77, 1).Value = 46040.5625
$detailed.Cells.Item(77, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(77, 2).Value = -8.452769999999999
$detailed.Cells.Item(77, 3).Value = "forecast"
$detailed.Cells.Item(77, 4).Value = 46040
$detailed.Cells.Item(77, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(77, 5).Value = "ON"

$detailed.Cells.Item(78, 1).Value = 46040.58333333334
$detailed.Cells.Item(78, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(78, 2).Value = -10
$detailed.Cells.Item(78, 3).Value = "forecast"
$detailed.Cells.Item(78, 4).Value = 46040
$detailed.Cells.Item(78, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(78, 5).Value = "ON"

$detailed.Cells.Item(79, 1).Value = 46040.60416666666
$detailed.Cells.Item(79, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(79, 2).Value = -14
$detailed.Cells.Item(79, 3).Value = "forecast"
$detailed.Cells.Item(79, 4).Value = 46040
$detailed.Cells.Item(79, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(79, 5).Value = "ON"

$detailed.Cells.Item(80, 1).Value = 46040.625
$detailed.Cells.Item(80, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(80, 2).Value = -19.76767
$detailed.Cells.Item(80, 3).Value = "forecast"
$detailed.Cells.Item(80, 4).Value = 46040
$detailed.Cells.Item(80, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(80, 5).Value = "ON"

$detailed.Cells.Item(81, 1).Value = 46040.64583333334
$detailed.Cells.Item(81, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(81, 2).Value = -20.57529
$detailed.Cells.Item(81, 3).Value = "forecast"
$detailed.Cells.Item(81, 4).Value = 46040
$detailed.Cells.Item(81, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(81, 5).Value = "ON"

$detailed.Cells.Item(82, 1).Value = 46040.66666666666
$detailed.Cells.Item(82, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(82, 2).Value = -14.35069
$detailed.Cells.Item(82, 3).Value = "forecast"
$detailed.Cells.Item(82, 4).Value = 46040
$detailed.Cells.Item(82, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(82, 5).Value = "ON"

$detailed.Cells.Item(83, 1).Value = 46040.6875
$detailed.Cells.Item(83, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(83, 2).Value = -14.23795
$detailed.Cells.Item(83, 3).Value = "forecast"
$detailed.Cells.Item(83, 4).Value = 46040
$detailed.Cells.Item(83, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(83, 5).Value = "ON"

$detailed.Cells.Item(84, 1).Value = 46040.70833333334
$detailed.Cells.Item(84, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(84, 2).Value = -12.30534
$detailed.Cells.Item(84, 3).Value = "forecast"
$detailed.Cells.Item(84, 4).Value = 46040
$detailed.Cells.Item(84, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(84, 5).Value = "ON"

$detailed.Cells.Item(85, 1).Value = 46040.72916666666
$detailed.Cells.Item(85, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(85, 2).Value = -6.59911
$detailed.Cells.Item(85, 3).Value = "forecast"
$detailed.Cells.Item(85, 4).Value = 46040
$detailed.Cells.Item(85, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(85, 5).Value = "ON"

$detailed.Cells.Item(86, 1).Value = 46040.75
$detailed.Cells.Item(86, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(86, 2).Value = 1.94062
$detailed.Cells.Item(86, 3).Value = "forecast"
$detailed.Cells.Item(86, 4).Value = 46040
$detailed.Cells.Item(86, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(86, 5).Value = "ON"

$detailed.Cells.Item(87, 1).Value = 46040.77083333334
$detailed.Cells.Item(87, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(87, 2).Value = 34.32564
$detailed.Cells.Item(87, 3).Value = "forecast"
$detailed.Cells.Item(87, 4).Value = 46040
$detailed.Cells.Item(87, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(87, 5).Value = "ON"

$detailed.Cells.Item(88, 1).Value = 46040.79166666666
$detailed.Cells.Item(88, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(88, 2).Value = 36.0601
$detailed.Cells.Item(88, 3).Value = "forecast"
$detailed.Cells.Item(88, 4).Value = 46040
$detailed.Cells.Item(88, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(88, 5).Value = "OFF"

$detailed.Cells.Item(89, 1).Value = 46040.8125
$detailed.Cells.Item(89, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(89, 2).Value = 36.2
$detailed.Cells.Item(89, 3).Value = "forecast"
$detailed.Cells.Item(89, 4).Value = 46040
$detailed.Cells.Item(89, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(89, 5).Value = "OFF"

$detailed.Cells.Item(90, 1).Value = 46040.83333333334
$detailed.Cells.Item(90, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(90, 2).Value = 43.32062
$detailed.Cells.Item(90, 3).Value = "forecast"
$detailed.Cells.Item(90, 4).Value = 46040
$detailed.Cells.Item(90, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(90, 5).Value = "OFF"

$detailed.Cells.Item(91, 1).Value = 46040.85416666666
$detailed.Cells.Item(91, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(91, 2).Value = 36.2
$detailed.Cells.Item(91, 3).Value = "forecast"
$detailed.Cells.Item(91, 4).Value = 46040
$detailed.Cells.Item(91, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(91, 5).Value = "OFF"

$detailed.Cells.Item(92, 1).Value = 46040.875
$detailed.Cells.Item(92, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(92, 2).Value = 36.2
$detailed.Cells.Item(92, 3).Value = "forecast"
$detailed.Cells.Item(92, 4).Value = 46040
$detailed.Cells.Item(92, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(92, 5).Value = "OFF"

$detailed.Cells.Item(93, 1).Value = 46040.89583333334
$detailed.Cells.Item(93, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(93, 2).Value = 36.0601
$detailed.Cells.Item(93, 3).Value = "forecast"
$detailed.Cells.Item(93, 4).Value = 46040
$detailed.Cells.Item(93, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(93, 5).Value = "OFF"

$detailed.Cells.Item(94, 1).Value = 46040.91666666666
$detailed.Cells.Item(94, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(94, 2).Value = 47.10825
$detailed.Cells.Item(94, 3).Value = "forecast"
$detailed.Cells.Item(94, 4).Value = 46040
$detailed.Cells.Item(94, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(94, 5).Value = "OFF"

$detailed.Cells.Item(95, 1).Value = 46040.9375
$detailed.Cells.Item(95, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(95, 2).Value = 36.2
$detailed.Cells.Item(95, 3).Value = "forecast"
$detailed.Cells.Item(95, 4).Value = 46040
$detailed.Cells.Item(95, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(95, 5).Value = "OFF"

$detailed.Cells.Item(96, 1).Value = 46040.95833333334
$detailed.Cells.Item(96, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(96, 2).Value = 36.0601
$detailed.Cells.Item(96, 3).Value = "forecast"
$detailed.Cells.Item(96, 4).Value = 46040
$detailed.Cells.Item(96, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(96, 5).Value = "OFF"

$detailed.Cells.Item(97, 1).Value = 46040.97916666666
$detailed.Cells.Item(97, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$detailed.Cells.Item(97, 2).Value = 36.06
$detailed.Cells.Item(97, 3).Value = "forecast"
$detailed.Cells.Item(97, 4).Value = 46040
$detailed.Cells.Item(97, 4).NumberFormat = "YYYY-MM-DD"
$detailed.Cells.Item(97, 5).Value = "OFF"

